$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.484.92"
$ws.Range("E2").Value = "  +5.71%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.053.99"
$ws.Range("E3").Value = "  +4.04%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.81"
$ws.Range("E5").Value = "  +3.49%  "

$ws.Range("E6").Value = "  +2.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "66.49"
$ws.Range("E7").Value = "  +17.37%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +7.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.84"
$ws.Range("E10").Value = "  +3.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0769"
$ws.Range("E11").Value = "  +5.23%  "

$ws.Range("E12").Value = "  +1.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.911"
$ws.Range("E13").Value = "  -3.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.95"
$ws.Range("E14").Value = "  +4.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.354.55"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.01"
$ws.Range("E16").Value = "  +25.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.60"
$ws.Range("E17").Value = "  +6.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.063.03"
$ws.Range("E18").Value = "  +4.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.269.93"
$ws.Range("E19").Value = "  +5.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.76"
$ws.Range("E20").Value = "  +3.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0879"
$ws.Range("E21").Value = "  +4.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.46"
$ws.Range("E22").Value = "  +6.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.50"
$ws.Range("E23").Value = "  +3.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.68"
$ws.Range("E24").Value = "  +4.94%  "

$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("E26").Value = "  +5.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.89"
$ws.Range("E27").Value = "  +9.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.65"
$ws.Range("E28").Value = "  -1.29%  "

$ws.Range("E29").Value = "  +4.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.118"
$ws.Range("E30").Value = "  +29.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.28"
$ws.Range("E31").Value = "  +9.26%  "

$ws.Range("E32").Value = "  +3.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.21"
$ws.Range("E33").Value = "  +8.87%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.71"
$ws.Range("E34").Value = "  +10.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0625"
$ws.Range("E35").Value = "  +6.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.47"
$ws.Range("E36").Value = "  +4.53%  "

$ws.Range("E37").Value = "  +4.47%  "

$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.04"
$ws.Range("E39").Value = "  +17.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.98"
$ws.Range("E40").Value = "  +33.22%  "

$ws.Range("E41").Value = "  +17.57%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.24"
$ws.Range("E42").Value = "  +2.72%  "

$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.01"
$ws.Range("E43").Value = "  +4.75%  "

$ws.Range("E44").Value = "  +6.52%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.21"
$ws.Range("E46").Value = "  +8.88%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "96.26"
$ws.Range("E47").Value = "  +5.98%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.99"
$ws.Range("E48").Value = "  +7.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.424.13"
$ws.Range("E49").Value = "  +3.67%  "

$ws.Range("E50").Value = "  +2.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.87"
$ws.Range("E51").Value = "  +2.09%  "
